$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.275.23'
$ws.Range("E2").Value = '  -2.28%  '
$ws.Range("D3").Value = '1.560.00'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '206.28'
$ws.Range("E5").Value = '  -3.15%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  -4.37%  '
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("E9").Value = '  -2.74%  '
$ws.Range("D10").Value = '17.74'
$ws.Range("E10").Value = '  -3.90%  '
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").Value = '1.775.77'
$ws.Range("E12").Value = '  -3.77%  '
$ws.Range("D13").Value = '1.562.21'
$ws.Range("E13").Value = '  -3.55%  '
$ws.Range("E14").Value = '  -3.71%  '
$ws.Range("D15").Value = '0.507'
$ws.Range("E15").Value = '  -3.16%  '
$ws.Range("D16").Value = '25.281.91'
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").Value = '59.22'
$ws.Range("E17").Value = '  -3.68%  '
$ws.Range("E18").Value = '  -3.59%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = '186.72'
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("D21").Value = '4.11'
$ws.Range("E21").Value = '  -3.05%  '
$ws.Range("D22").Value = '9.25'
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("E23").Value = '  -2.85%  '
$ws.Range("D24").Value = '0.130'
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("D25").Value = '141.03'
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("D28").Value = '14.91'
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("E29").Value = '  -4.35%  '
$ws.Range("E30").Value = '  -6.87%  '
$ws.Range("D31").Value = '0.0465'
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  -3.90%  '
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("E35").Value = '  -4.17%  '
$ws.Range("D36").Value = '1.088.21'
$ws.Range("E36").Value = '  -3.39%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '2.35'
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("B38").Value = 'PaxDollar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("D39").Value = '0.496'
$ws.Range("E39").Value = '  -3.41%  '
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("D41").Value = '0.771'
$ws.Range("E41").Value = '  -8.26%  '
$ws.Range("D42").Value = '0.799'
$ws.Range("E42").Value = '  +6.64%  '
$ws.Range("D43").Value = '93.09'
$ws.Range("E43").Value = '  -5.12%  '
$ws.Range("D44").Value = '5.10'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").Value = '1.695.78'
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").Value = '1.48'
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D48").Value = '52.50'
$ws.Range("E48").Value = '  -3.03%  '
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.404'
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  -0.36%  '
